$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 200; this shifts existing rows 200-273 down to 201-274
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(200).Insert()

# Populate the new row 200 with the new weekly record.
$ws.Cells.Item(200, 1).Value2  = 10
$ws.Cells.Item(200, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(200, 3).Value2  = "La Araucanía"
$ws.Cells.Item(200, 4).Value2  = 44559
$ws.Cells.Item(200, 5).Value2  = 9
$ws.Cells.Item(200, 6).Value2  = "Fruta"
$ws.Cells.Item(200, 7).Value2  = 100108
$ws.Cells.Item(200, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(200, 9).Value2  = 100108002
$ws.Cells.Item(200, 10).Value2 = "Mango"
$ws.Cells.Item(200, 11).Value2 = "Sin especificar"
$ws.Cells.Item(200, 12).Value2 = "Primera"
$ws.Cells.Item(200, 13).Value2 = 190
$ws.Cells.Item(200, 14).Value2 = 7000
$ws.Cells.Item(200, 15).Value2 = 7000
$ws.Cells.Item(200, 16).Value2 = 7000
$ws.Cells.Item(200, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(200, 18).Value2 = "Brasil"
$ws.Cells.Item(200, 19).Value2 = 1750
$ws.Cells.Item(200, 20).Value2 = 4
